# Scheduled market-price refresh: updates currentAveragePrice* / Leve cost-profit
# columns (H:N) for the leves whose item prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: Enchanted Silver Ink
$ws.Range("H28").Value = 436.6316
$ws.Range("I28").Value = 470.2857
$ws.Range("J28").Value = 342.4
$ws.Range("K28").Value = 470.2857
$ws.Range("L28").Value = 342.4
$ws.Range("M28").Value = 14.71429999999998
$ws.Range("N28").Value = -1312.4

# Row 41: Enchanted Mythril Ink
$ws.Range("H41").Value = 3203.1428
$ws.Range("I41").Value = 371.88235
$ws.Range("J41").Value = 5877.1113
$ws.Range("K41").Value = 371.88235
$ws.Range("L41").Value = 5877.1113
$ws.Range("M41").Value = 68.11765000000003
$ws.Range("N41").Value = -6757.1113

# Row 62: Enchanted Mythrite Ink
$ws.Range("H62").Value = 7276.385
$ws.Range("I62").Value = 7603.7812
$ws.Range("J62").Value = 5779.7144
$ws.Range("K62").Value = 7603.7812
$ws.Range("L62").Value = 5779.7144
$ws.Range("M62").Value = -6979.7812
$ws.Range("N62").Value = -7027.7144

# Row 65: Enchanted Mythrite Ink
$ws.Range("H65").Value = 7276.385
$ws.Range("I65").Value = 7603.7812
$ws.Range("J65").Value = 5779.7144
$ws.Range("K65").Value = 38018.906
$ws.Range("L65").Value = 28898.572
$ws.Range("M65").Value = -34898.906
$ws.Range("N65").Value = -35138.572

# Row 98: Enchanted Durium Ink
$ws.Range("H98").Value = 1793.1305
$ws.Range("I98").Value = 1006.13336
$ws.Range("K98").Value = 1006.13336
$ws.Range("M98").Value = 491.86664

# Row 116: Growth Formula Kappa
$ws.Range("H116").Value = 32292.555
$ws.Range("I116").Value = 53259.285
$ws.Range("J116").Value = 2939.1333
$ws.Range("K116").Value = 53259.285
$ws.Range("L116").Value = 2939.1333
$ws.Range("M116").Value = -49817.285
$ws.Range("N116").Value = -9823.1333

# Row 122: Enchanted High Durium Ink
$ws.Range("H122").Value = 1793.1305
$ws.Range("I122").Value = 1006.13336
$ws.Range("K122").Value = 3018.40008
$ws.Range("M122").Value = -568.4000800000003

# Row 124: Luncheon Toadskin Codex
$ws.Range("H124").Value = 29737.5
$ws.Range("J124").Value = 29737.5
$ws.Range("L124").Value = 29737.5
$ws.Range("N124").Value = -39557.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 1942.1383
$ws.Range("I31").Value = 1503.1666
$ws.Range("J31").Value = 2716.7942
$ws.Range("K31").Value = 1503.1666
$ws.Range("L31").Value = 2716.7942
$ws.Range("M31").Value = -1208.1666
$ws.Range("N31").Value = -3306.7942

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 1942.1383
$ws.Range("I34").Value = 1503.1666
$ws.Range("J34").Value = 2716.7942
$ws.Range("K34").Value = 1503.1666
$ws.Range("L34").Value = 2716.7942
$ws.Range("M34").Value = -1301.1666
$ws.Range("N34").Value = -3120.7942

# Row 55: Mythril Lance
$ws.Range("H55").Value = 13400
$ws.Range("I55").Value = 10500
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 10500
$ws.Range("L55").Value = 25000
$ws.Range("M55").Value = -10185
$ws.Range("N55").Value = -25630

# Row 58: Mahogany Lumber
$ws.Range("H58").Value = 3292.6206
$ws.Range("J58").Value = 2746.913
$ws.Range("L58").Value = 2746.913
$ws.Range("N58").Value = -3152.913

# Row 99: Pine Lumber
$ws.Range("H99").Value = 57805.61
$ws.Range("I99").Value = 113280.89
$ws.Range("J99").Value = 2330.3333
$ws.Range("K99").Value = 113280.89
$ws.Range("L99").Value = 2330.3333
$ws.Range("M99").Value = -111782.89
$ws.Range("N99").Value = -5326.3333

# Row 126: Red Pine Lumber
$ws.Range("H126").Value = 57805.61
$ws.Range("I126").Value = 113280.89
$ws.Range("J126").Value = 2330.3333
$ws.Range("K126").Value = 339842.67
$ws.Range("L126").Value = 6990.999899999999
$ws.Range("M126").Value = -337372.67
$ws.Range("N126").Value = -11930.9999

# Row 136: Dark Mahogany Lumber
$ws.Range("H136").Value = 3292.6206
$ws.Range("J136").Value = 2746.913
$ws.Range("L136").Value = 8240.739
$ws.Range("N136").Value = -13340.739

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Lavender Oil
$ws.Range("H23").Value = 483.66666
$ws.Range("I23").Value = 53
$ws.Range("J23").Value = 627.2222
$ws.Range("K23").Value = 159
$ws.Range("L23").Value = 1881.6666
$ws.Range("M23").Value = 76
$ws.Range("N23").Value = -2351.6666

# Row 113: Night Vinegar
$ws.Range("H113").Value = 9615935
$ws.Range("I113").Value = 14286253
$ws.Range("J113").Value = 574.64703
$ws.Range("K113").Value = 42858759
$ws.Range("L113").Value = 1723.94109
$ws.Range("M113").Value = -42856589
$ws.Range("N113").Value = -6063.94109

$ws = $wb.Worksheets.Item("GSM")
# Row 47: Peridot Choker
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# Row 55: Peridot Earrings
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# Row 136: Pink Beryl
$ws.Range("H136").Value = 24331.5
$ws.Range("J136").Value = 24331.5
$ws.Range("L136").Value = 72994.5
$ws.Range("N136").Value = -78094.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Aldgoat Leather
$ws.Range("H22").Value = 501.625
$ws.Range("I22").Value = 516.1429000000001
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 516.1429000000001
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -221.1429000000001
$ws.Range("N22").Value = -990

# Row 27: Aldgoat Leather
$ws.Range("H27").Value = 501.625
$ws.Range("I27").Value = 516.1429000000001
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 516.1429000000001
$ws.Range("L27").Value = 400
$ws.Range("M27").Value = -409.1429000000001
$ws.Range("N27").Value = -614

# Row 48: Fingerless Peisteskin Gloves
$ws.Range("H48").Value = 15561.833
$ws.Range("I48").Value = 3490
$ws.Range("J48").Value = 21597.75
$ws.Range("K48").Value = 3490
$ws.Range("L48").Value = 21597.75
$ws.Range("M48").Value = -2829
$ws.Range("N48").Value = -22919.75

# Row 55: Peiste Leather
$ws.Range("H55").Value = 214.17647
$ws.Range("I55").Value = 195
$ws.Range("J55").Value = 241.57143
$ws.Range("K55").Value = 195
$ws.Range("L55").Value = 241.57143
$ws.Range("M55").Value = -22
$ws.Range("N55").Value = -587.57143

$ws = $wb.Worksheets.Item("WVR")
# Row 126: Snow Linen
$ws.Range("H126").Value = 926.6667
$ws.Range("I126").Value = 507.69232
$ws.Range("J126").Value = 1199
$ws.Range("K126").Value = 1523.07696
$ws.Range("L126").Value = 3597
$ws.Range("M126").Value = 946.9230400000001
$ws.Range("N126").Value = -8537

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 3060.4348
$ws.Range("I132").Value = 2241.5715
$ws.Range("J132").Value = 4334.222
$ws.Range("K132").Value = 6724.7145
$ws.Range("L132").Value = 13002.666
$ws.Range("M132").Value = -4194.7145
$ws.Range("N132").Value = -18062.666

# Row 137: Sarcenet Slops of Aiming
$ws.Range("H137").Value = 47175.25
$ws.Range("J137").Value = 47175.25
$ws.Range("L137").Value = 47175.25
$ws.Range("N137").Value = -57375.25
